$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 15:40"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 4500130
$ws.Range("C4").Value = 1787
$ws.Range("E4").Value = 2158150
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = 152388

# --- Row 16: Arabia Saudita ---
$ws.Range("B16").Value = 272590
$ws.Range("C16").Value = 1759
$ws.Range("D16").Value = 228569
$ws.Range("E16").Value = 41205
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 2816

# --- Row 24: Irak ---
$ws.Range("B24").Value = 118300
$ws.Range("C24").Value = 2968
$ws.Range("D24").Value = 83461
$ws.Range("E24").Value = 30236
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 4603

# --- Row 44: Paises Bajos ---
$ws.Range("B44").Value = 53621
$ws.Range("C44").Value = 247
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 6147

# --- Row 64: Uzbekistan ---
$ws.Range("B64").Value = 22374
$ws.Range("C64").Value = 481
$ws.Range("D64").Value = 12668
$ws.Range("E64").Value = 9577

# --- Row 73: Chequia ---
$ws.Range("C73").Value = 275

# --- Rows 80/81: Bosnia y Herzegovina <-> Estado de Palestina swap rank,
# with Estado de Palestina's numbers updated for the day ---
$ws.Range("A80").Value = "Estado de Palestina"
$ws.Range("B80").Value = 11284
$ws.Range("C80").Value = 346
$ws.Range("D80").Value = 4833
$ws.Range("E80").Value = 6372
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 79

$ws.Range("A81").Value = "Bosnia y Herzegovina"
$ws.Range("B81").Value = 11127
$ws.Range("C81").Value = 361
$ws.Range("D81").Value = 5441
$ws.Range("E81").Value = 5370
$ws.Range("G81").Value = 19
$ws.Range("H81").Value = 316

# --- Row 86: Noruega ---
$ws.Range("B86").Value = 9158
$ws.Range("C86").Value = 8
$ws.Range("E86").Value = 151

# --- Row 90: Finlandia ---
$ws.Range("D90").Value = 6950
$ws.Range("E90").Value = 135

# --- Row 98: Republica de Yibuti ---
$ws.Range("B98").Value = 5081
$ws.Range("C98").Value = 13
$ws.Range("D98").Value = 4999
$ws.Range("E98").Value = 24

# --- Row 142: Liberia ---
$ws.Range("B142").Value = 1179
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 664
$ws.Range("E142").Value = 443

# --- Rows 160/161: Tanzania <-> Lesoto swap rank,
# with Lesoto's numbers updated for the day ---
$ws.Range("A160").Value = "Lesoto"
$ws.Range("B160").Value = 576
$ws.Range("C160").Value = 71
$ws.Range("D160").Value = 141
$ws.Range("E160").Value = 422
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 13

$ws.Range("A161").Value = "Tanzania"
$ws.Range("B161").Value = 509
$ws.Range("D161").Value = 183
$ws.Range("E161").Value = 305
$ws.Range("H161").Value = 21

# --- Row 181: Trinidad yTobago ---
$ws.Range("B181").Value = 154
$ws.Range("C181").Value = 1
$ws.Range("E181").Value = 18
